# Update worksheet "Dados": every cell in column F that contains the
# label "Outro" should now contain "Outra" instead.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    if ($cell.Value2 -eq "Outro") {
        $cell.Value = "Outra"
    }
}
